# Lecture partielle de l'EDT M1 MIAGE.
# Update the weekday labels and the corresponding date serials on the
# "Liste" sheet: the schedule read-out moved from the 2023 edition of the
# calendar to the 2026 one (a +1096 day shift), which also changes which
# weekday each class date falls on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates (column A) - same day-of-month/month, shifted 3 years (+1096 days)
$ws.Range("A2").Value = 46064.0
$ws.Range("A5").Value = 46073.0
$ws.Range("A9").Value = 46091.0
$ws.Range("A11").Value = 46092.0
$ws.Range("A14").Value = 46105.0
$ws.Range("A16").Value = 46108.0

# Weekday labels (column B) matching the new dates above
$ws.Range("B2").Value = "mercredi"
$ws.Range("B5").Value = "vendredi"
$ws.Range("B9").Value = "mardi"
$ws.Range("B11").Value = "mercredi"
$ws.Range("B14").Value = "mardi"
$ws.Range("B16").Value = "vendredi"
